$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in new data row describing the DB preparation export, matching the
# existing "footer" rows of the checklist table.
$ws.Range("A20").Value = "Подготовка на БД"
$ws.Range("D20").Value = "*"

# Update the selected cell to mirror where the user clicked last.
$ws.Range("F19").Select()
